$wb = $excel.ActiveWorkbook

# The "展览" (sheet 1) and "全部类型" (sheet 4) worksheets both hold the same
# 12-event convention table (header row + rows 2-13). This update drops the
# oldest event (which used to live in row 2) and pulls every other event's
# details up by one row, while keeping the running index in column A (1..11)
# untouched. After the shift, a handful of "interest count" (column F)
# figures were also refreshed to newer numbers, and the now-duplicate last
# row is removed so the sheet shrinks from 13 rows to 12.

$targetSheets = @(1, 4)

foreach ($idx in $targetSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # Shift columns B:I (everything except the running index in column A)
    # up by one data row: row 3's content becomes row 2's content, etc.
    $src = $ws.Range("B3:I13")
    $dst = $ws.Range("B2")
    $src.Copy($dst)

    # Drop the now-duplicated trailing row so the table goes back to 12 rows.
    $ws.Rows.Item(13).Delete()

    # A few "interest count" values (column F) were bumped beyond a plain
    # shift-up once the new data came in.
    $ws.Range("F2").Value2 = 6432
    $ws.Range("F4").Value2 = 24
    $ws.Range("F6").Value2 = 1924
    $ws.Range("F7").Value2 = 1465
    $ws.Range("F8").Value2 = 301
    $ws.Range("F9").Value2 = 987
    $ws.Range("F10").Value2 = 317
    $ws.Range("F11").Value2 = 5603
}
